# The document originally has two paragraphs:
#   1) (empty paragraph)
#   2) "31233" + <bookmarkStart/End name="_GoBack"/> + "hehheh"
#
# Target edit (per the diff):
#   - remove the "_GoBack" bookmark from paragraph 2 (leaving "31233hehheh"
#     as plain, unbookmarked runs)
#   - append a brand-new paragraph 3 containing the text "gdfdg", with the
#     "_GoBack" bookmark now wrapping the end of that new paragraph instead

$d = $word.ActiveDocument

# --- Step 1: drop the existing "_GoBack" bookmark -------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: add the new paragraph with text "gdfdg" and a relocated ------
# --- "_GoBack" bookmark at the end of that paragraph's text ---------------
$end = $d.Content
$end.Collapse(0)   # wdCollapseEnd -> collapse to the very end of the document

$xml = '<?xml version="1.0"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' +
                        '<w:p>' +
                            '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
                            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>gdfdg</w:t></w:r>' +
                            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                            '<w:bookmarkEnd w:id="0"/>' +
                        '</w:p>' +
                    '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

$end.InsertXML($xml)
